# Auto-generated: update price/profit columns (H-N) across Sheets tables
# per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7893.773
$ws.Range("I62").Value = 9446.9375
$ws.Range("K62").Value = 9446.9375
$ws.Range("M62").Value = -8822.9375

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 7893.773
$ws.Range("I65").Value = 9446.9375
$ws.Range("K65").Value = 47234.6875
$ws.Range("M65").Value = -44114.6875

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1500
$ws.Range("I111").Value = 1500
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 4500
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = -1433
$ws.Range("N111").Value = -10634

# ALC row 117
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 10508.5
$ws.Range("I129").Value = 500
$ws.Range("J129").Value = 13844.667
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 41534.001
$ws.Range("M129").Value = 3500
$ws.Range("N129").Value = -51534.001

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5480.5654
$ws.Range("I137").Value = 951.13336
$ws.Range("K137").Value = 2853.40008
$ws.Range("M137").Value = -303.4000800000003

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1176.375
$ws.Range("I2").Value = 942.2
$ws.Range("J2").Value = 1566.6666
$ws.Range("K2").Value = 942.2
$ws.Range("L2").Value = 1566.6666
$ws.Range("M2").Value = -829.2
$ws.Range("N2").Value = -1792.6666

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2160.4167
$ws.Range("I61").Value = 2160.4167
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2160.4167
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1948.4167
$ws.Range("N61").ClearContents()

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1777.5161
$ws.Range("I74").Value = 1994.65
$ws.Range("J74").Value = 1382.7273
$ws.Range("K74").Value = 1994.65
$ws.Range("L74").Value = 1382.7273
$ws.Range("M74").Value = -1120.65
$ws.Range("N74").Value = -3130.7273

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1777.5161
$ws.Range("I77").Value = 1994.65
$ws.Range("J77").Value = 1382.7273
$ws.Range("K77").Value = 9973.25
$ws.Range("L77").Value = 6913.636500000001
$ws.Range("M77").Value = -5605.25
$ws.Range("N77").Value = -15649.6365

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1727.8889
$ws.Range("I102").Value = 1522.5
$ws.Range("J102").Value = 1892.2
$ws.Range("K102").Value = 1522.5
$ws.Range("L102").Value = 1892.2
$ws.Range("M102").Value = 99.5
$ws.Range("N102").Value = -5136.2

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1176.375
$ws.Range("I116").Value = 942.2
$ws.Range("J116").Value = 1566.6666
$ws.Range("K116").Value = 942.2
$ws.Range("L116").Value = 1566.6666
$ws.Range("M116").Value = 1351.8
$ws.Range("N116").Value = -6154.6666

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 71887.56
$ws.Range("I132").Value = 86266.914
$ws.Range("K132").Value = 258800.742
$ws.Range("M132").Value = -256270.742

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2160.4167
$ws.Range("I136").Value = 2160.4167
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6481.250100000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3931.250100000001
$ws.Range("N136").ClearContents()

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1176.375
$ws.Range("I3").Value = 942.2
$ws.Range("J3").Value = 1566.6666
$ws.Range("K3").Value = 942.2
$ws.Range("L3").Value = 1566.6666
$ws.Range("M3").Value = -828.2
$ws.Range("N3").Value = -1794.6666

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2614.7837
$ws.Range("I134").Value = 1628.6786
$ws.Range("J134").Value = 5682.6665
$ws.Range("K134").Value = 4886.0358
$ws.Range("L134").Value = 17047.9995
$ws.Range("M134").Value = -2351.0358
$ws.Range("N134").Value = -22117.9995

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 49067.355
$ws.Range("I31").Value = 58098.668
$ws.Range("J31").Value = 21973.428
$ws.Range("K31").Value = 58098.668
$ws.Range("L31").Value = 21973.428
$ws.Range("M31").Value = -57803.668
$ws.Range("N31").Value = -22563.428

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 49067.355
$ws.Range("I34").Value = 58098.668
$ws.Range("J34").Value = 21973.428
$ws.Range("K34").Value = 58098.668
$ws.Range("L34").Value = 21973.428
$ws.Range("M34").Value = -57896.668
$ws.Range("N34").Value = -22377.428

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2550.111
$ws.Range("I58").Value = 901.53845
$ws.Range("J58").Value = 6836.4
$ws.Range("K58").Value = 901.53845
$ws.Range("L58").Value = 6836.4
$ws.Range("M58").Value = -698.53845
$ws.Range("N58").Value = -7242.4

# CRP row 92
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 30800.5
$ws.Range("J92").Value = 30800.5
$ws.Range("L92").Value = 30800.5
$ws.Range("N92").Value = -35792.5

# CRP row 114
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2313.4138
$ws.Range("J132").Value = 3930.6
$ws.Range("L132").Value = 11791.8
$ws.Range("N132").Value = -16851.8

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 10205499
$ws.Range("I134").Value = 1203.6342
$ws.Range("K134").Value = 3610.9026
$ws.Range("M134").Value = -1075.9026

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2550.111
$ws.Range("I136").Value = 901.53845
$ws.Range("J136").Value = 6836.4
$ws.Range("K136").Value = 2704.61535
$ws.Range("L136").Value = 20509.2
$ws.Range("M136").Value = -154.61535
$ws.Range("N136").Value = -25609.2

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 100005
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 100005
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 300015
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -300239

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 100005
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 100005
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 900045
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -905115

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 5298049
$ws.Range("I140").Value = 2012
$ws.Range("J140").Value = 10112628
$ws.Range("K140").Value = 6036
$ws.Range("L140").Value = 30337884
$ws.Range("M140").Value = -856
$ws.Range("N140").Value = -30348244

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 30078.861
$ws.Range("I132").Value = 1582
$ws.Range("J132").Value = 74859.64
$ws.Range("K132").Value = 4746
$ws.Range("L132").Value = 224578.92
$ws.Range("M132").Value = -2216
$ws.Range("N132").Value = -229638.92

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 28279.36
$ws.Range("I132").Value = 42382.16
$ws.Range("J132").Value = 3095.7856
$ws.Range("K132").Value = 127146.48
$ws.Range("L132").Value = 9287.356800000001
$ws.Range("M132").Value = -124616.48
$ws.Range("N132").Value = -14347.3568

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1665.24
$ws.Range("I136").Value = 1021.4
$ws.Range("J136").Value = 2631
$ws.Range("K136").Value = 3064.2
$ws.Range("L136").Value = 7893
$ws.Range("M136").Value = -514.1999999999998
$ws.Range("N136").Value = -12993

# WVR row 80
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 41255.668
$ws.Range("J80").Value = 41255.668
$ws.Range("L80").Value = 41255.668
$ws.Range("N80").Value = -43251.668

# WVR row 83
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 41255.668
$ws.Range("J83").Value = 41255.668
$ws.Range("L83").Value = 123767.004
$ws.Range("N83").Value = -133751.004

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1545.2046
$ws.Range("I132").Value = 1268.8387
$ws.Range("J132").Value = 2204.2307
$ws.Range("K132").Value = 3806.5161
$ws.Range("L132").Value = 6612.6921
$ws.Range("M132").Value = -1276.5161
$ws.Range("N132").Value = -11672.6921

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2395823.8
$ws.Range("I136").Value = 2551935.8
$ws.Range("J136").Value = 1667301.5
$ws.Range("K136").Value = 7655807.399999999
$ws.Range("L136").Value = 5001904.5
$ws.Range("M136").Value = -7653257.399999999
$ws.Range("N136").Value = -5007004.5

